# feat: add 2022-Q4 data
#
# 1. Insert a new worksheet "2022-Q4" right after "总计" (and before "2022-Q3"),
#    populated with the Q4 fund-holding detail rows.
# 2. Insert a new summary row at the top of "总计" for 2022-Q4, shifting the
#    existing quarter rows down by one and renumbering the running index in
#    column A.

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)          # "总计" - always sheet 1
$q3Before = $wb.Worksheets.Item(2)          # current "2022-Q3", used only to anchor the Add()

# ---------------------------------------------------------------------------
# 1) Build the new "2022-Q4" sheet, inserted before the current "2022-Q3" tab.
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($q3Before, $null)
$q4.Name = "2022-Q4"

# NOTE: inserting a sheet shifts every worksheet index after it, so any
# worksheet reference captured *before* the Add() call becomes unreliable for
# formatting operations (Copy/PasteSpecial silently drops formats through a
# stale reference, even though plain value reads still work). Re-resolve the
# "2022-Q3" sheet by name now that it exists in the workbook again.
$q3 = $wb.Worksheets.Item("2022-Q3")

# Header row, copying the bold/bordered/centered style used by every other
# quarter sheet's header row (B1:H1) from the existing "2022-Q3" sheet.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q4.Cells.Item(1, $i + 2).Value = $headers[$i]
}
$q3.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

# Data rows: column A is a styled running index (numeric, same style as the
# rest of the workbook's index column); B/C are text; D:G look numeric but
# are stored as text (matching the source data); H is a real number.
$q4.Cells.Item(2, 1).Value = 0
$q4.Cells.Item(3, 1).Value = 1
$q3.Range("A2:A3").Copy()
$q4.Range("A2:A3").PasteSpecial(-4122)

$q4Rows = @(
    @("580006", "东吴新经济混合A", "1.16", "89.15", "4.03", "0.0467", 9),
    @("012617", "东吴新经济混合C", "0.42", "89.15", "4.03", "0.0169", 9)
)

for ($r = 0; $r -lt $q4Rows.Length; $r++) {
    $row = $q4Rows[$r]
    $excelRow = $r + 2

    $q4.Cells.Item($excelRow, 2).NumberFormat = "@"
    $q4.Cells.Item($excelRow, 2).Value = $row[0]
    $q4.Cells.Item($excelRow, 2).Style = "Normal"

    $q4.Cells.Item($excelRow, 3).NumberFormat = "@"
    $q4.Cells.Item($excelRow, 3).Value = $row[1]
    $q4.Cells.Item($excelRow, 3).Style = "Normal"

    $q4.Cells.Item($excelRow, 4).NumberFormat = "@"
    $q4.Cells.Item($excelRow, 4).Value = $row[2]
    $q4.Cells.Item($excelRow, 4).Style = "Normal"

    $q4.Cells.Item($excelRow, 5).NumberFormat = "@"
    $q4.Cells.Item($excelRow, 5).Value = $row[3]
    $q4.Cells.Item($excelRow, 5).Style = "Normal"

    $q4.Cells.Item($excelRow, 6).NumberFormat = "@"
    $q4.Cells.Item($excelRow, 6).Value = $row[4]
    $q4.Cells.Item($excelRow, 6).Style = "Normal"

    $q4.Cells.Item($excelRow, 7).NumberFormat = "@"
    $q4.Cells.Item($excelRow, 7).Value = $row[5]
    $q4.Cells.Item($excelRow, 7).Style = "Normal"

    $q4.Cells.Item($excelRow, 8).Value = $row[6]
}

# ---------------------------------------------------------------------------
# 2) Insert the new 2022-Q4 row at the top of the "总计" summary sheet.
# ---------------------------------------------------------------------------

# Capture the existing quarter rows (rows 2-5) before shifting them.
$existing = @()
for ($r = 2; $r -le 5; $r++) {
    $label = $summary.Cells.Item($r, 2).Value()
    $count = $summary.Cells.Item($r, 3).Value()
    $value = $summary.Cells.Item($r, 4).Value()
    $existing += , @($label, $count, $value)
}

# Re-write shifted down one row (now rows 3-6).
for ($i = $existing.Length - 1; $i -ge 0; $i--) {
    $destRow = $i + 3
    $row = $existing[$i]
    $summary.Cells.Item($destRow, 1).Value = $i + 1
    $summary.Cells.Item($destRow, 2).Value = $row[0]
    $summary.Cells.Item($destRow, 3).Value = $row[1]
    $summary.Cells.Item($destRow, 4).Value = $row[2]
}

# The bottom-most destination row (row 6) is brand new - give column A the
# same style as the rest of the index column.
$summary.Cells.Item(6, 1).Value = 4
$summary.Range("A5").Copy()
$summary.Range("A6").PasteSpecial(-4122)
$summary.Cells.Item(6, 1).Value = 4

# New row 2: 2022-Q4 summary values.
$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 2
$summary.Cells.Item(2, 4).Value = 0.06
